$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($row, $col, $text) {
    $cell = $ws.Cells.Item($row, $col)
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = "Normal"
}

# Row 2 - Bitcoin
Set-TextCell 2 4 "57.325.25"
Set-TextCell 2 5 "  +3.28%  "

# Row 3 - Ethereum
Set-TextCell 3 4 "3.068.91"
Set-TextCell 3 5 "  +5.57%  "

# Row 4 - TetherUSD
Set-TextCell 4 5 "  -0.17%  "

# Row 5 - BNB
Set-TextCell 5 4 "515.07"
Set-TextCell 5 5 "  +1.80%  "

# Row 6 - Solana
Set-TextCell 6 4 "141.38"
Set-TextCell 6 5 "  +6.69%  "

# Row 7 - USDC
Set-TextCell 7 5 "  -0.07%  "

# Row 8 - XRP
Set-TextCell 8 5 "  +3.61%  "

# Row 9 - Toncoin
Set-TextCell 9 4 "7.25"
Set-TextCell 9 5 "  +1.17%  "

# Row 10 - Dogecoin
Set-TextCell 10 5 "  +4.40%  "

# Row 11 - Cardano
Set-TextCell 11 4 "0.372"
Set-TextCell 11 5 "  +6.90%  "

# Row 12 - WrappedliquidstakedEther2.0
Set-TextCell 12 4 "3.593.34"
Set-TextCell 12 5 "  +5.41%  "

# Row 13 - TRON
Set-TextCell 13 5 "  +2.70%  "

# Row 14 - Avalanche
Set-TextCell 14 4 "25.49"
Set-TextCell 14 5 "  -0.90%  "

# Row 15 - ShibaInu
Set-TextCell 15 5 "  +4.11%  "

# Row 16 - WrappedBTC
Set-TextCell 16 4 "57.320.70"
Set-TextCell 16 5 "  +3.20%  "

# Row 17 - WrappedEther
Set-TextCell 17 4 "3.065.74"
Set-TextCell 17 5 "  +5.31%  "

# Row 18 - Polkadot
Set-TextCell 18 4 "5.94"
Set-TextCell 18 5 "  -0.73%  "

# Row 19 - Chainlink
Set-TextCell 19 5 "  +4.26%  "

# Row 20 - Uniswap
Set-TextCell 20 5 "  +6.39%  "

# Row 21 - BitcoinCash
Set-TextCell 21 4 "337.69"
Set-TextCell 21 5 "  +7.70%  "

# Row 22 - Dai
Set-TextCell 22 5 "  +0.30%  "

# Row 23 - Polygon
Set-TextCell 23 4 "0.499"
Set-TextCell 23 5 "  +3.69%  "

# Row 24 - Litecoin
Set-TextCell 24 4 "65.41"
Set-TextCell 24 5 "  +4.52%  "

# Row 25 - Kaspa
Set-TextCell 25 5 "  +6.81%  "

# Row 26 - Binance-PegBSC-USD
Set-TextCell 26 5 "  -0.03%  "

# Row 27 - PEPE
Set-TextCell 27 4 "0.0₃0944"
Set-TextCell 27 5 "  +12.93%  "

# Row 28 - RenderToken
Set-TextCell 28 4 "6.43"
Set-TextCell 28 5 "  +1.75%  "

# Row 29 - InternetComputer(DFINITY)
Set-TextCell 29 4 "7.07"
Set-TextCell 29 5 "  +2.68%  "

# Row 30 - PancakeSwap
Set-TextCell 30 5 "  +2.27%  "

# Row 31 - EthereumClassic
Set-TextCell 31 4 "20.75"
Set-TextCell 31 5 "  +5.50%  "

# Row 32 - Fetch.AI
Set-TextCell 32 5 "  +4.96%  "

# Row 33 - Monero
Set-TextCell 33 4 "154.46"
Set-TextCell 33 5 "  +3.67%  "

# Row 34 - NEARProtocol
Set-TextCell 34 5 "  +4.30%  "

# Row 35 - Aptos
Set-TextCell 35 4 "5.88"
Set-TextCell 35 5 "  +5.70%  "

# Row 36 - EnergySwap
Set-TextCell 36 4 "26.00"
Set-TextCell 36 5 "  +6.67%  "

# Row 37 - ImmutableX
Set-TextCell 37 5 "  +5.45%  "

# Row 38 - Hedera
Set-TextCell 38 4 "0.0672"
Set-TextCell 38 5 "  +4.64%  "

# Row 39 - RenzoRestakedETH
Set-TextCell 39 4 "3.104.78"
Set-TextCell 39 5 "  +5.60%  "

# Row 40 - OKB
Set-TextCell 40 4 "36.96"
Set-TextCell 40 5 "  +1.93%  "

# Row 41 - Mantle
Set-TextCell 41 4 "0.669"
Set-TextCell 41 5 "  +5.54%  "

# Row 42 - Filecoin
Set-TextCell 42 5 "  +4.40%  "

# Row 43 - FirstDigitalUSD
Set-TextCell 43 5 "  -0.19%  "

# Row 44 - Maker
Set-TextCell 44 4 "2.249.05"
Set-TextCell 44 5 "  +7.08%  "

# Row 45 - VeChain
Set-TextCell 45 4 "0.0252"
Set-TextCell 45 5 "  +8.99%  "

# Row 46 - Stacks
Set-TextCell 46 5 "  +4.80%  "

# Row 47 - ONDO
Set-TextCell 47 4 "0.951"
Set-TextCell 47 5 "  +4.83%  "

# Row 48 - InjectiveProtocol
Set-TextCell 48 4 "20.06"
Set-TextCell 48 5 "  +8.19%  "

# Row 49 - Cosmos
Set-TextCell 49 4 "5.84"
Set-TextCell 49 5 "  -0.59%  "

# Row 50 - Stellar
Set-TextCell 50 5 "  +4.08%  "

# Row 51 - SuiNetwork -> dogwifhat
Set-TextCell 51 2 "dogwifhat"
Set-TextCell 51 3 "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
Set-TextCell 51 4 "1.73"
Set-TextCell 51 5 "  +4.15%  "
